$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update expiry timestamps for the first two sample cookies
$ws.Range("G2").Value = 1736828073
$ws.Range("G3").Value = 1738815273

# Third sample cookie is now a session cookie (no expiry) -> expiry column holds text "session"
$ws.Range("G4").Value = "session"
